# Update the cryptos list (Price column D, Volume(1h) column E) with the
# latest scraped values. D-column values are written with a leading
# apostrophe so Excel stores them as literal text (matching the original
# inlineStr cells) instead of auto-converting numeric-looking strings
# (e.g. "0.211", "613.19") into real numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + "69.422.60"
$ws.Range("E2").Value = "  -1.92%  "

$ws.Range("D3").Value = "'" + "3.487.15"
$ws.Range("E3").Value = "  -2.14%  "

$ws.Range("D4").Value = "'" + "0.998"
$ws.Range("E4").Value = "  -0.25%  "

$ws.Range("D5").Value = "'" + "613.19"
$ws.Range("E5").Value = "  +5.43%  "

$ws.Range("D6").Value = "'" + "189.08"
$ws.Range("E6").Value = "  +0.95%  "

$ws.Range("D7").Value = "'" + "0.626"
$ws.Range("E7").Value = "  -0.23%  "

$ws.Range("E8").Value = "  -0.09%  "

$ws.Range("D9").Value = "'" + "0.211"

$ws.Range("D10").Value = "'" + "0.648"
$ws.Range("E10").Value = "  -0.47%  "

$ws.Range("D11").Value = "'" + "52.75"
$ws.Range("E11").Value = "  -3.21%  "

$ws.Range("D12").Value = "'" + "0.0000306"
$ws.Range("E12").Value = "  -3.89%  "

$ws.Range("D13").Value = "'" + "9.45"
$ws.Range("E13").Value = "  -0.76%  "

$ws.Range("D14").Value = "'" + "4.042.81"
$ws.Range("E14").Value = "  -2.25%  "

$ws.Range("D15").Value = "'" + "613.15"
$ws.Range("E15").Value = "  +7.03%  "

$ws.Range("D16").Value = "'" + "69.461.16"
$ws.Range("E16").Value = "  -1.97%  "

$ws.Range("D17").Value = "'" + "18.88"
$ws.Range("E17").Value = "  -1.69%  "

$ws.Range("D18").Value = "'" + "12.52"
$ws.Range("E18").Value = "  -2.50%  "

$ws.Range("D19").Value = "'" + "3.481.64"
$ws.Range("E19").Value = "  -2.85%  "

$ws.Range("E20").Value = "  -0.36%  "

$ws.Range("D21").Value = "'" + "0.982"
$ws.Range("E21").Value = "  -2.16%  "

$ws.Range("E22").Value = "  -3.27%  "

$ws.Range("D23").Value = "'" + "106.01"
$ws.Range("E23").Value = "  +12.42%  "

$ws.Range("D24").Value = "'" + "4.70"
$ws.Range("E24").Value = "  +2.59%  "

$ws.Range("D25").Value = "'" + "5.10"
$ws.Range("E25").Value = "  +4.21%  "

$ws.Range("D26").Value = "'" + "3.01"
$ws.Range("E26").Value = "  +1.44%  "

$ws.Range("D27").Value = "'" + "10.93"
$ws.Range("E27").Value = "  -2.33%  "

$ws.Range("D28").Value = "'" + "9.68"
$ws.Range("E28").Value = "  +3.51%  "

$ws.Range("D29").Value = "'" + "33.51"
$ws.Range("E29").Value = "  +2.30%  "

$ws.Range("D30").Value = "'" + "6.90"
$ws.Range("E30").Value = "  -4.18%  "

$ws.Range("D31").Value = "'" + "12.55"
$ws.Range("E31").Value = "  +1.91%  "

$ws.Range("D32").Value = "'" + "3.94"
$ws.Range("E32").Value = "  +3.57%  "

$ws.Range("E33").Value = "  -1.73%  "

$ws.Range("D34").Value = "'" + "63.28"
$ws.Range("E34").Value = "  -0.11%  "

$ws.Range("E35").Value = "  -5.36%  "

$ws.Range("D36").Value = "'" + "0.999"
$ws.Range("E36").Value = "  -0.07%  "

$ws.Range("D37").Value = "'" + "3.615.96"
$ws.Range("E37").Value = "  -0.53%  "

$ws.Range("D38").Value = "'" + "3.63"
$ws.Range("E38").Value = "  +4.98%  "

$ws.Range("E39").Value = "  -4.51%  "

$ws.Range("D40").Value = "'" + "506.45"
$ws.Range("E40").Value = "  -5.55%  "

$ws.Range("D41").Value = "'" + "36.46"
$ws.Range("E41").Value = "  -4.29%  "

$ws.Range("D42").Value = "'" + "0.0" + [char]0x2083 + "0768"
$ws.Range("E42").Value = "  -5.34%  "

$ws.Range("E43").Value = "  -3.49%  "

$ws.Range("D44").Value = "'" + "0.0459"
$ws.Range("E44").Value = "  -2.67%  "

$ws.Range("E45").Value = "  -1.89%  "

$ws.Range("E46").Value = "  +1.99%  "

$ws.Range("E47").Value = "  -4.43%  "

$ws.Range("E48").Value = "  +0.35%  "

$ws.Range("D49").Value = "'" + "8.70"
$ws.Range("E49").Value = "  -6.93%  "

$ws.Range("D50").Value = "'" + "131.01"
$ws.Range("E50").Value = "  -4.14%  "

$ws.Range("E51").Value = "  -7.54%  "
